$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maçlar")

# Add the new match result links for the 30.06.2025 records
# (written in this order so the new shared strings land at the same
# indices as the reference edit: 36 = rbFAYTWC6z4, 37 = mjBKimYNaCk)
$ws.Range("H15").Value = "https://youtu.be/rbFAYTWC6z4"
$ws.Range("H14").Value = "https://youtu.be/mjBKimYNaCk"

# Update the active selection to reflect the last edited cell
$ws.Range("H14").Select()
